$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('LP1912')
$ws1.Cells.Item(2,1).Value = 'Última actualización: 17:39:57'
$ws1.Cells.Item(3,1).Value = 'Total filas: 416'
$ws1.Cells.Item(82,1).Value = '06:52:23'
$ws1.Cells.Item(82,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(82,4).Value = 110
$ws1.Cells.Item(83,1).Value = '08:39:56'
$ws1.Cells.Item(83,3).Value = '14_ABASTO'
$ws1.Cells.Item(83,4).Value = 3
$ws1.Cells.Item(108,1).Value = '08:21:27'
$ws1.Cells.Item(108,3).Value = '17_ROMERO'
$ws1.Cells.Item(108,4).Value = 61
$ws1.Cells.Item(109,1).Value = '07:46:15'
$ws1.Cells.Item(109,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(109,4).Value = 96
$ws1.Cells.Item(120,1).Value = '08:50:00'
$ws1.Cells.Item(120,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(120,4).Value = 45
$ws1.Cells.Item(121,1).Value = '08:57:11'
$ws1.Cells.Item(121,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(121,4).Value = 38
$ws1.Cells.Item(159,1).Value = '10:57:58'
$ws1.Cells.Item(159,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(159,4).Value = 7
$ws1.Cells.Item(160,1).Value = '10:28:12'
$ws1.Cells.Item(160,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(160,4).Value = 36
$ws1.Cells.Item(171,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(172,3).Value = '17_ROMERO'
$ws1.Cells.Item(212,3).Value = '215A_EL PATO'
$ws1.Cells.Item(213,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(224,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(226,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(251,1).Value = '11:51:05'
$ws1.Cells.Item(251,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(251,4).Value = 90
$ws1.Cells.Item(252,1).Value = '12:44:21'
$ws1.Cells.Item(252,3).Value = '10_OLMOS'
$ws1.Cells.Item(252,4).Value = 37
$ws1.Cells.Item(264,3).Value = '215A_EL PATO'
$ws1.Cells.Item(265,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(303,1).Value = '14:16:51'
$ws1.Cells.Item(303,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(303,4).Value = 57
$ws1.Cells.Item(304,1).Value = '14:40:41'
$ws1.Cells.Item(304,3).Value = '10_OLMOS'
$ws1.Cells.Item(304,4).Value = 33
$ws1.Cells.Item(354,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(355,3).Value = '225_GOMEZ'
$ws1.Cells.Item(386,1).Value = '17:39:57'
$ws1.Cells.Item(386,2).Value = '17:39'
$ws1.Cells.Item(386,4).Value = 0
$ws1.Cells.Item(387,1).Value = '17:39:57'
$ws1.Cells.Item(387,2).Value = '17:39'
$ws1.Cells.Item(387,4).Value = 0
$ws1.Cells.Item(388,1).Value = '17:39:57'
$ws1.Cells.Item(388,2).Value = '17:39'
$ws1.Cells.Item(388,3).Value = '17_ROMERO'
$ws1.Cells.Item(388,4).Value = 0
$ws1.Cells.Item(389,1).Value = '15:51:40'
$ws1.Cells.Item(389,2).Value = '17:40'
$ws1.Cells.Item(389,3).Value = '215B_EL PATO'
$ws1.Cells.Item(389,4).Value = 109
$ws1.Cells.Item(390,1).Value = '16:52:27'
$ws1.Cells.Item(390,2).Value = '17:40'
$ws1.Cells.Item(390,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(390,4).Value = 48
$ws1.Cells.Item(391,1).Value = '16:45:22'
$ws1.Cells.Item(391,2).Value = '17:41'
$ws1.Cells.Item(391,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(391,4).Value = 56
$ws1.Cells.Item(392,1).Value = '17:39:57'
$ws1.Cells.Item(392,2).Value = '17:44'
$ws1.Cells.Item(392,3).Value = '15_ABASTO'
$ws1.Cells.Item(392,4).Value = 5
$ws1.Cells.Item(393,1).Value = '16:45:22'
$ws1.Cells.Item(393,2).Value = '17:45'
$ws1.Cells.Item(393,3).Value = '15_ABASTO'
$ws1.Cells.Item(393,4).Value = 60
$ws1.Cells.Item(394,1).Value = '15:51:40'
$ws1.Cells.Item(394,2).Value = '17:50'
$ws1.Cells.Item(394,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(394,4).Value = 119
$ws1.Cells.Item(395,1).Value = '17:39:57'
$ws1.Cells.Item(395,2).Value = '17:51'
$ws1.Cells.Item(395,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(395,4).Value = 12
$ws1.Cells.Item(396,1).Value = '16:14:52'
$ws1.Cells.Item(396,2).Value = '17:52'
$ws1.Cells.Item(396,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(396,4).Value = 98
$ws1.Cells.Item(397,1).Value = '17:39:57'
$ws1.Cells.Item(397,2).Value = '17:52'
$ws1.Cells.Item(397,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(397,4).Value = 13
$ws1.Cells.Item(398,1).Value = '17:14:54'
$ws1.Cells.Item(398,2).Value = '17:59'
$ws1.Cells.Item(398,3).Value = '10_OLMOS'
$ws1.Cells.Item(398,4).Value = 45
$ws1.Cells.Item(399,1).Value = '17:39:57'
$ws1.Cells.Item(399,2).Value = '18:00'
$ws1.Cells.Item(399,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(399,4).Value = 21
$ws1.Cells.Item(400,1).Value = '17:39:57'
$ws1.Cells.Item(400,2).Value = '18:03'
$ws1.Cells.Item(400,3).Value = '17_ROMERO'
$ws1.Cells.Item(400,4).Value = 24
$ws1.Cells.Item(401,1).Value = '17:14:54'
$ws1.Cells.Item(401,2).Value = '18:04'
$ws1.Cells.Item(401,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(401,4).Value = 50
$ws1.Cells.Item(402,1).Value = '16:14:52'
$ws1.Cells.Item(402,2).Value = '18:04'
$ws1.Cells.Item(402,3).Value = '17_ROMERO'
$ws1.Cells.Item(402,4).Value = 110
$ws1.Cells.Item(403,2).Value = '18:08'
$ws1.Cells.Item(403,3).Value = '14_ABASTO'
$ws1.Cells.Item(403,4).Value = 76
$ws1.Cells.Item(404,1).Value = '17:39:57'
$ws1.Cells.Item(404,2).Value = '18:15'
$ws1.Cells.Item(404,3).Value = '15_ABASTO'
$ws1.Cells.Item(404,4).Value = 36
$ws1.Cells.Item(405,1).Value = '17:39:57'
$ws1.Cells.Item(405,2).Value = '18:15'
$ws1.Cells.Item(405,3).Value = '10_OLMOS'
$ws1.Cells.Item(405,4).Value = 36
$ws1.Cells.Item(406,1).Value = '17:39:57'
$ws1.Cells.Item(406,2).Value = '18:20'
$ws1.Cells.Item(406,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(406,4).Value = 66
$ws1.Cells.Item(407,1).Value = '16:32:38'
$ws1.Cells.Item(407,2).Value = '18:21'
$ws1.Cells.Item(407,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(407,4).Value = 109
$ws1.Cells.Item(407,5).Value = 'LP1912'
$ws1.Cells.Item(408,1).Value = '17:39:57'
$ws1.Cells.Item(408,2).Value = '18:24'
$ws1.Cells.Item(408,3).Value = '14_ABASTO'
$ws1.Cells.Item(408,4).Value = 45
$ws1.Cells.Item(408,5).Value = 'LP1912'
$ws1.Cells.Item(409,1).Value = '16:32:38'
$ws1.Cells.Item(409,2).Value = '18:27'
$ws1.Cells.Item(409,3).Value = '215C_EL PATO'
$ws1.Cells.Item(409,4).Value = 115
$ws1.Cells.Item(409,5).Value = 'LP1912'
$ws1.Cells.Item(410,1).Value = '16:45:22'
$ws1.Cells.Item(410,2).Value = '18:28'
$ws1.Cells.Item(410,3).Value = '215C_EL PATO'
$ws1.Cells.Item(410,4).Value = 103
$ws1.Cells.Item(410,5).Value = 'LP1912'
$ws1.Cells.Item(411,1).Value = '17:14:54'
$ws1.Cells.Item(411,2).Value = '18:31'
$ws1.Cells.Item(411,3).Value = '11X44_ETCHEVERRY'
$ws1.Cells.Item(411,4).Value = 77
$ws1.Cells.Item(411,5).Value = 'LP1912'
$ws1.Cells.Item(412,1).Value = '16:45:22'
$ws1.Cells.Item(412,2).Value = '18:32'
$ws1.Cells.Item(412,3).Value = '11X44_ETCHEVERRY'
$ws1.Cells.Item(412,4).Value = 107
$ws1.Cells.Item(412,5).Value = 'LP1912'
$ws1.Cells.Item(413,1).Value = '17:39:57'
$ws1.Cells.Item(413,2).Value = '18:36'
$ws1.Cells.Item(413,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(413,4).Value = 57
$ws1.Cells.Item(413,5).Value = 'LP1912'
$ws1.Cells.Item(414,1).Value = '17:14:54'
$ws1.Cells.Item(414,2).Value = '18:47'
$ws1.Cells.Item(414,3).Value = '14X44_ABASTO'
$ws1.Cells.Item(414,4).Value = 93
$ws1.Cells.Item(414,5).Value = 'LP1912'
$ws1.Cells.Item(415,1).Value = '16:52:27'
$ws1.Cells.Item(415,2).Value = '18:48'
$ws1.Cells.Item(415,3).Value = '14X44_ABASTO'
$ws1.Cells.Item(415,4).Value = 116
$ws1.Cells.Item(415,5).Value = 'LP1912'
$ws1.Cells.Item(416,1).Value = '17:14:54'
$ws1.Cells.Item(416,2).Value = '18:58'
$ws1.Cells.Item(416,3).Value = '215A_EL PATO'
$ws1.Cells.Item(416,4).Value = 104
$ws1.Cells.Item(416,5).Value = 'LP1912'
$ws1.Cells.Item(417,1).Value = '17:14:54'
$ws1.Cells.Item(417,2).Value = '19:04'
$ws1.Cells.Item(417,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(417,4).Value = 110
$ws1.Cells.Item(417,5).Value = 'LP1912'
$ws1.Cells.Item(418,1).Value = '17:14:54'
$ws1.Cells.Item(418,2).Value = '19:10'
$ws1.Cells.Item(418,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(418,4).Value = 116
$ws1.Cells.Item(418,5).Value = 'LP1912'
$ws1.Cells.Item(419,1).Value = '17:39:57'
$ws1.Cells.Item(419,2).Value = '19:16'
$ws1.Cells.Item(419,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(419,4).Value = 97
$ws1.Cells.Item(419,5).Value = 'LP1912'
$ws1.Cells.Item(420,1).Value = '17:39:57'
$ws1.Cells.Item(420,2).Value = '19:20'
$ws1.Cells.Item(420,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(420,4).Value = 101
$ws1.Cells.Item(420,5).Value = 'LP1912'
$ws1.Cells.Item(421,1).Value = '17:39:57'
$ws1.Cells.Item(421,2).Value = '19:29'
$ws1.Cells.Item(421,3).Value = '225_GOMEZ'
$ws1.Cells.Item(421,4).Value = 110
$ws1.Cells.Item(421,5).Value = 'LP1912'

$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws2.Cells.Item(2,1).Value = 'Última actualización: 17:39:57'
$ws2.Cells.Item(3,1).Value = 'Total filas: 42'
$ws2.Cells.Item(43,1).Value = '17:39:57'
$ws2.Cells.Item(43,2).Value = '17:39'
$ws2.Cells.Item(43,4).Value = 0
$ws2.Cells.Item(44,1).Value = '15:51:40'
$ws2.Cells.Item(44,2).Value = '17:40'
$ws2.Cells.Item(44,3).Value = '215B_EL PATO'
$ws2.Cells.Item(44,4).Value = 109
$ws2.Cells.Item(45,1).Value = '16:32:38'
$ws2.Cells.Item(45,2).Value = '18:27'
$ws2.Cells.Item(45,4).Value = 115
$ws2.Cells.Item(46,1).Value = '16:45:22'
$ws2.Cells.Item(46,2).Value = '18:28'
$ws2.Cells.Item(46,3).Value = '215C_EL PATO'
$ws2.Cells.Item(46,4).Value = 103
$ws2.Cells.Item(47,1).Value = '17:14:54'
$ws2.Cells.Item(47,2).Value = '18:58'
$ws2.Cells.Item(47,3).Value = '215A_EL PATO'
$ws2.Cells.Item(47,4).Value = 104
$ws2.Cells.Item(47,5).Value = 'LP1912'

$ws3 = $wb.Worksheets.Item('6203-6173')
$ws3.Cells.Item(2,1).Value = 'Última actualización: 17:39:57'

